$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuel specs")

# Insert a new row at 42 (shifts existing rows 42-94 down to 43-95)
# and populate it with the "Renewable Diesel" entry.
$ws.Rows(42).Insert()
$ws.Range("A42").Value = "Renewable Diesel"
$ws.Range("B42").Value = 40.669004600898425
$ws.Range("C42").Value = 793.37856525946859

# Make "Fuel specs" the active sheet / tab, with a specific scroll
# position and active selection, matching the saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A43").Select()
